$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 775.8946999999999
$ws.Range("J17").Value = 759.7406999999999
$ws.Range("L17").Value = 2279.2221
$ws.Range("N17").Value = -2615.2221
$ws.Range("H32").Value = 4159.5947
$ws.Range("J32").Value = 4421.364
$ws.Range("L32").Value = 4421.364
$ws.Range("N32").Value = -5073.364
$ws.Range("H40").Value = 3668.2058
$ws.Range("I40").Value = 1398.75
$ws.Range("K40").Value = 1398.75
$ws.Range("M40").Value = -1223.75
$ws.Range("H58").Value = 412.33334
$ws.Range("I58").Value = 94.8
$ws.Range("K58").Value = 284.4
$ws.Range("M58").Value = -134.4
$ws.Range("H76").Value = 6496.385
$ws.Range("I76").Value = 6334.75
$ws.Range("J76").Value = 6755
$ws.Range("K76").Value = 6334.75
$ws.Range("L76").Value = 6755
$ws.Range("M76").Value = -6019.75
$ws.Range("N76").Value = -7385
$ws.Range("H79").Value = 6496.385
$ws.Range("I79").Value = 6334.75
$ws.Range("J79").Value = 6755
$ws.Range("K79").Value = 6334.75
$ws.Range("L79").Value = 6755
$ws.Range("M79").Value = -5242.75
$ws.Range("N79").Value = -8939
$ws.Range("H100").Value = 2905.3635
$ws.Range("I100").Value = 3017.7778
$ws.Range("J100").Value = 2399.5
$ws.Range("K100").Value = 3017.7778
$ws.Range("L100").Value = 2399.5
$ws.Range("M100").Value = -2476.7778
$ws.Range("N100").Value = -3481.5
$ws.Range("H112").Value = 4758.9
$ws.Range("J112").Value = 4758.9
$ws.Range("L112").Value = 14276.7
$ws.Range("N112").Value = -16492.7
$ws.Range("H131").Value = 4563.0625
$ws.Range("I131").Value = 4268.6665
$ws.Range("J131").Value = 5446.25
$ws.Range("K131").Value = 12805.9995
$ws.Range("L131").Value = 16338.75
$ws.Range("M131").Value = -7765.999500000002
$ws.Range("N131").Value = -26418.75
$ws.Range("H138").Value = 1974.38
$ws.Range("I138").Value = 841.9524
$ws.Range("J138").Value = 2794.4138
$ws.Range("K138").Value = 2525.8572
$ws.Range("L138").Value = 8383.241399999999
$ws.Range("M138").Value = 2614.1428
$ws.Range("N138").Value = -18663.2414

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H74").Value = 168901.27
$ws.Range("I74").Value = 72526
$ws.Range("J74").Value = 441964.5
$ws.Range("K74").Value = 72526
$ws.Range("L74").Value = 441964.5
$ws.Range("M74").Value = -71652
$ws.Range("N74").Value = -443712.5
$ws.Range("H77").Value = 168901.27
$ws.Range("I77").Value = 72526
$ws.Range("J77").Value = 441964.5
$ws.Range("K77").Value = 362630
$ws.Range("L77").Value = 2209822.5
$ws.Range("M77").Value = -358262
$ws.Range("N77").Value = -2218558.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H105").Value = 3127295.2
$ws.Range("I105").Value = 3908644.5
$ws.Range("J105").Value = 1897.5
$ws.Range("K105").Value = 3908644.5
$ws.Range("L105").Value = 1897.5
$ws.Range("M105").Value = -3906897.5
$ws.Range("N105").Value = -5391.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 7999
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 7999
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = 7999
$ws.Range("N5").Value = -8223
$ws.Range("L5").ClearContents()
$ws.Range("H31").Value = 31801.355
$ws.Range("I31").Value = 1602.6
$ws.Range("J31").Value = 60112.688
$ws.Range("K31").Value = 1602.6
$ws.Range("L31").Value = 60112.688
$ws.Range("M31").Value = -1307.6
$ws.Range("N31").Value = -60702.688
$ws.Range("H34").Value = 31801.355
$ws.Range("I34").Value = 1602.6
$ws.Range("J34").Value = 60112.688
$ws.Range("K34").Value = 1602.6
$ws.Range("L34").Value = 60112.688
$ws.Range("M34").Value = -1400.6
$ws.Range("N34").Value = -60516.688
$ws.Range("H86").Value = 14759.839
$ws.Range("I86").Value = 13745.842
$ws.Range("J86").Value = 16365.333
$ws.Range("K86").Value = 13745.842
$ws.Range("L86").Value = 16365.333
$ws.Range("M86").Value = -12622.842
$ws.Range("N86").Value = -18611.333
$ws.Range("H89").Value = 14759.839
$ws.Range("I89").Value = 13745.842
$ws.Range("J89").Value = 16365.333
$ws.Range("K89").Value = 68729.21000000001
$ws.Range("L89").Value = 81826.66500000001
$ws.Range("M89").Value = -63113.21000000001
$ws.Range("N89").Value = -93058.66500000001
$ws.Range("H94").Value = 1066.1875
$ws.Range("I94").Value = 736.6667
$ws.Range("J94").Value = 1263.9
$ws.Range("K94").Value = 736.6667
$ws.Range("L94").Value = 1263.9
$ws.Range("M94").Value = -285.6667
$ws.Range("N94").Value = -2165.9
$ws.Range("H127").Value = 70000
$ws.Range("J127").Value = 70000
$ws.Range("L127").Value = 70000
$ws.Range("N127").Value = -79920
$ws.Range("H132").Value = 91936.09
$ws.Range("I132").Value = 60398.53
$ws.Range("J132").Value = 225970.75
$ws.Range("K132").Value = 181195.59
$ws.Range("L132").Value = 677912.25
$ws.Range("M132").Value = -178665.59
$ws.Range("N132").Value = -682972.25
$ws.Range("H134").Value = 30169.188
$ws.Range("I134").Value = 33764.965
$ws.Range("K134").Value = 101294.895
$ws.Range("M134").Value = -98759.89499999999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 55755.5
$ws.Range("I12").Value = 177967.6
$ws.Range("K12").Value = 533902.8
$ws.Range("M12").Value = -533729.8
$ws.Range("H116").Value = 7545.1816
$ws.Range("I116").Value = 100
$ws.Range("J116").Value = 8289.700000000001
$ws.Range("K116").Value = 300
$ws.Range("L116").Value = 24869.1
$ws.Range("M116").Value = 3142
$ws.Range("N116").Value = -31753.1
$ws.Range("H131").Value = 16668978
$ws.Range("I131").Value = 10417506
$ws.Range("J131").Value = 19610846
$ws.Range("K131").Value = 31252518
$ws.Range("L131").Value = 58832538
$ws.Range("M131").Value = -31247478
$ws.Range("N131").Value = -58842618

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 84508.336
$ws.Range("J51").Value = 84508.336
$ws.Range("L51").Value = 84508.336
$ws.Range("N51").Value = -85526.336
$ws.Range("H70").Value = 7697842.5
$ws.Range("I70").Value = 10531264
$ws.Range("K70").Value = 10531264
$ws.Range("M70").Value = -10530994
$ws.Range("H73").Value = 7697842.5
$ws.Range("I73").Value = 10531264
$ws.Range("K73").Value = 10531264
$ws.Range("M73").Value = -10530328
$ws.Range("H122").Value = 214584.14
$ws.Range("I122").Value = 279793
$ws.Range("J122").Value = 5915.8
$ws.Range("K122").Value = 839379
$ws.Range("L122").Value = 17747.4
$ws.Range("M122").Value = -836929
$ws.Range("N122").Value = -22647.4
$ws.Range("H126").Value = 7157914
$ws.Range("I126").Value = 3499532.2
$ws.Range("K126").Value = 10498596.6
$ws.Range("M126").Value = -10496126.6

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3669.4211
$ws.Range("I7").Value = 2203.3572
$ws.Range("J7").Value = 7774.4
$ws.Range("K7").Value = 2203.3572
$ws.Range("L7").Value = 7774.4
$ws.Range("M7").Value = -2091.3572
$ws.Range("N7").Value = -7998.4
$ws.Range("H16").Value = 1474.9474
$ws.Range("I16").Value = 1247.9166
$ws.Range("J16").Value = 1864.1428
$ws.Range("K16").Value = 1247.9166
$ws.Range("L16").Value = 1864.1428
$ws.Range("M16").Value = -1077.9166
$ws.Range("N16").Value = -2204.1428
$ws.Range("H122").Value = 4897.4443
$ws.Range("J122").Value = 8927.6
$ws.Range("L122").Value = 26782.8
$ws.Range("N122").Value = -31682.8
$ws.Range("H126").Value = 3669.4211
$ws.Range("I126").Value = 2203.3572
$ws.Range("J126").Value = 7774.4
$ws.Range("K126").Value = 6610.071599999999
$ws.Range("L126").Value = 23323.2
$ws.Range("M126").Value = -4140.071599999999
$ws.Range("N126").Value = -28263.2
$ws.Range("H133").Value = 141999
$ws.Range("J133").Value = 141999
$ws.Range("L133").Value = 141999
$ws.Range("N133").Value = -147059
$ws.Range("H134").Value = 44000
$ws.Range("J134").Value = 44000
$ws.Range("L134").Value = 44000
$ws.Range("N134").Value = -54140
$ws.Range("H136").Value = 50034.836
$ws.Range("I136").Value = 77535.44500000001
$ws.Range("K136").Value = 232606.335
$ws.Range("M136").Value = -230056.335

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 50801
$ws.Range("J103").Value = 50801
$ws.Range("L103").Value = 50801
$ws.Range("N103").Value = -53145
$ws.Range("H106").Value = 36299.5
$ws.Range("J106").Value = 36299.5
$ws.Range("L106").Value = 36299.5
$ws.Range("N106").Value = -38823.5
$ws.Range("H107").Value = 55559144
$ws.Range("I107").Value = 62503830
$ws.Range("K107").Value = 187511490
$ws.Range("M107").Value = -187509570
